$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) -> F7, F11, F20 updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F7").Value = 9642
$ws1.Range("F11").Value = 2760
$ws1.Range("F20").Value = 1359

# Sheet "全部类型" (All types) -> F8, F12, F21 updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F8").Value = 9642
$ws4.Range("F12").Value = 2760
$ws4.Range("F21").Value = 1359
